# chore(runtime): publish files + archive (2025-11-22 15:05:04)
#
# Adds two newly-completed KHL matches (uids 897800, 897799) played on
# 2025-11-21 to Matches_SOG, then rolls their shots-on-goal totals forward
# into the Shots_HA / Shots_Summary aggregate sheets and bumps the
# Meta_ext "as of" timestamp + build_version.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell while forcing it to stay TEXT (Excel
# otherwise auto-converts a purely-numeric-looking string, like the uid
# "897800", into a Number). We force the "@" text format just long enough
# to take the assignment, then restore the default "Normal" style so no
# stray formatting is left behind on the cell.
# ---------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) Matches_SOG: append the two new match rows (502, 503)
# ---------------------------------------------------------------------
$matches = $wb.Worksheets.Item("Matches_SOG")

Set-TextValue $matches.Cells.Item(502, 1) "897800"
Set-TextValue $matches.Cells.Item(502, 2) "2025-11-21T17:00:00"
$matches.Cells.Item(502, 3).Value = "Барыс"
$matches.Cells.Item(502, 4).Value = "Салават Юлаев"
$matches.Cells.Item(502, 5).Value = 31
$matches.Cells.Item(502, 6).Value = 27
$matches.Cells.Item(502, 7).Value = "khl_text"

Set-TextValue $matches.Cells.Item(503, 1) "897799"
Set-TextValue $matches.Cells.Item(503, 2) "2025-11-21T19:30:00"
$matches.Cells.Item(503, 3).Value = "Нефтехимик"
$matches.Cells.Item(503, 4).Value = "Ак Барс"
$matches.Cells.Item(503, 5).Value = 27
$matches.Cells.Item(503, 6).Value = 31
$matches.Cells.Item(503, 7).Value = "khl_text"

# ---------------------------------------------------------------------
# 2) Shots_HA: refresh as_of_utc for every team row, then roll the two
#    new matches' shots into the home/away aggregate totals for the four
#    teams involved (Барыс, Салават Юлаев, Нефтехимик, Ак Барс).
# ---------------------------------------------------------------------
$shotsHA = $wb.Worksheets.Item("Shots_HA")

for ($r = 2; $r -le 23; $r++) {
    $shotsHA.Cells.Item($r, 4).Value = "2025-11-21T19:30:00Z"
}

# Барыс (row 7) played at home vs Салават Юлаев: sog_home=31, sog_away=27
$row = 7
$gpHome = $shotsHA.Cells.Item($row, 5).Value2 + 1
$hogfTotal = $shotsHA.Cells.Item($row, 7).Value2 + 31
$hogaTotal = $shotsHA.Cells.Item($row, 8).Value2 + 27
$shotsHA.Cells.Item($row, 5).Value = $gpHome
$shotsHA.Cells.Item($row, 7).Value = $hogfTotal
$shotsHA.Cells.Item($row, 8).Value = $hogaTotal
$shotsHA.Cells.Item($row, 9).Value = [math]::Round($hogfTotal / $gpHome, 1)
$shotsHA.Cells.Item($row, 10).Value = [math]::Round($hogaTotal / $gpHome, 1)

# Салават Юлаев (row 16) played away at Барыс: sog_away=27, sog_home=31
$row = 16
$gpAway = $shotsHA.Cells.Item($row, 6).Value2 + 1
$aogfTotal = $shotsHA.Cells.Item($row, 11).Value2 + 27
$aogaTotal = $shotsHA.Cells.Item($row, 12).Value2 + 31
$shotsHA.Cells.Item($row, 6).Value = $gpAway
$shotsHA.Cells.Item($row, 11).Value = $aogfTotal
$shotsHA.Cells.Item($row, 12).Value = $aogaTotal
$shotsHA.Cells.Item($row, 13).Value = [math]::Round($aogfTotal / $gpAway, 1)
$shotsHA.Cells.Item($row, 14).Value = [math]::Round($aogaTotal / $gpAway, 1)

# Нефтехимик (row 14) played at home vs Ак Барс: sog_home=27, sog_away=31
$row = 14
$gpHome = $shotsHA.Cells.Item($row, 5).Value2 + 1
$hogfTotal = $shotsHA.Cells.Item($row, 7).Value2 + 27
$hogaTotal = $shotsHA.Cells.Item($row, 8).Value2 + 31
$shotsHA.Cells.Item($row, 5).Value = $gpHome
$shotsHA.Cells.Item($row, 7).Value = $hogfTotal
$shotsHA.Cells.Item($row, 8).Value = $hogaTotal
$shotsHA.Cells.Item($row, 9).Value = [math]::Round($hogfTotal / $gpHome, 1)
$shotsHA.Cells.Item($row, 10).Value = [math]::Round($hogaTotal / $gpHome, 1)

# Ак Барс (row 5) played away at Нефтехимик: sog_away=31, sog_home=27
$row = 5
$gpAway = $shotsHA.Cells.Item($row, 6).Value2 + 1
$aogfTotal = $shotsHA.Cells.Item($row, 11).Value2 + 31
$aogaTotal = $shotsHA.Cells.Item($row, 12).Value2 + 27
$shotsHA.Cells.Item($row, 6).Value = $gpAway
$shotsHA.Cells.Item($row, 11).Value = $aogfTotal
$shotsHA.Cells.Item($row, 12).Value = $aogaTotal
$shotsHA.Cells.Item($row, 13).Value = [math]::Round($aogfTotal / $gpAway, 1)
$shotsHA.Cells.Item($row, 14).Value = [math]::Round($aogaTotal / $gpAway, 1)

# ---------------------------------------------------------------------
# 3) Shots_Summary: refresh as_of_utc for every team row, then roll the
#    combined (home+away) shots-for/against totals for the same teams.
# ---------------------------------------------------------------------
$shotsSummary = $wb.Worksheets.Item("Shots_Summary")

for ($r = 2; $r -le 23; $r++) {
    $shotsSummary.Cells.Item($r, 4).Value = "2025-11-21T19:30:00Z"
}

# Барыс (row 7): GP_total+1, SOG_total+=31, SOGA_total+=27
$row = 7
$gp = $shotsSummary.Cells.Item($row, 5).Value2 + 1
$sog = $shotsSummary.Cells.Item($row, 6).Value2 + 31
$soga = $shotsSummary.Cells.Item($row, 7).Value2 + 27
$shotsSummary.Cells.Item($row, 5).Value = $gp
$shotsSummary.Cells.Item($row, 6).Value = $sog
$shotsSummary.Cells.Item($row, 7).Value = $soga
$shotsSummary.Cells.Item($row, 8).Value = [math]::Round($sog / $gp, 1)
$shotsSummary.Cells.Item($row, 9).Value = [math]::Round($soga / $gp, 1)

# Салават Юлаев (row 16): GP_total+1, SOG_total+=27, SOGA_total+=31
$row = 16
$gp = $shotsSummary.Cells.Item($row, 5).Value2 + 1
$sog = $shotsSummary.Cells.Item($row, 6).Value2 + 27
$soga = $shotsSummary.Cells.Item($row, 7).Value2 + 31
$shotsSummary.Cells.Item($row, 5).Value = $gp
$shotsSummary.Cells.Item($row, 6).Value = $sog
$shotsSummary.Cells.Item($row, 7).Value = $soga
$shotsSummary.Cells.Item($row, 8).Value = [math]::Round($sog / $gp, 1)
$shotsSummary.Cells.Item($row, 9).Value = [math]::Round($soga / $gp, 1)

# Нефтехимик (row 14): GP_total+1, SOG_total+=27, SOGA_total+=31
$row = 14
$gp = $shotsSummary.Cells.Item($row, 5).Value2 + 1
$sog = $shotsSummary.Cells.Item($row, 6).Value2 + 27
$soga = $shotsSummary.Cells.Item($row, 7).Value2 + 31
$shotsSummary.Cells.Item($row, 5).Value = $gp
$shotsSummary.Cells.Item($row, 6).Value = $sog
$shotsSummary.Cells.Item($row, 7).Value = $soga
$shotsSummary.Cells.Item($row, 8).Value = [math]::Round($sog / $gp, 1)
$shotsSummary.Cells.Item($row, 9).Value = [math]::Round($soga / $gp, 1)

# Ак Барс (row 5): GP_total+1, SOG_total+=31, SOGA_total+=27
$row = 5
$gp = $shotsSummary.Cells.Item($row, 5).Value2 + 1
$sog = $shotsSummary.Cells.Item($row, 6).Value2 + 31
$soga = $shotsSummary.Cells.Item($row, 7).Value2 + 27
$shotsSummary.Cells.Item($row, 5).Value = $gp
$shotsSummary.Cells.Item($row, 6).Value = $sog
$shotsSummary.Cells.Item($row, 7).Value = $soga
$shotsSummary.Cells.Item($row, 8).Value = [math]::Round($sog / $gp, 1)
$shotsSummary.Cells.Item($row, 9).Value = [math]::Round($soga / $gp, 1)

# ---------------------------------------------------------------------
# 4) Meta_ext: bump as_of_utc + build_version
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Meta_ext")
$meta.Cells.Item(2, 2).Value = "2025-11-21T19:30:00Z"
$meta.Cells.Item(2, 4).Value = 84
